$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as Text (the source data stores
# these as inline strings, e.g. "66.723.26" or "0.137" -- left alone, Excel
# auto-converts numeric-looking text into a real number/float on assignment).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2" "66.720.20"
$ws.Range("E2").Value = "  -0.32%  "
Set-TextValue "D3" "3.778.51"
$ws.Range("E3").Value = "  -2.59%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "437.49"
$ws.Range("E5").Value = "  +2.16%  "
Set-TextValue "D6" "142.09"
$ws.Range("E6").Value = "  +7.59%  "
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  -8.85%  "
Set-TextValue "D11" "0.0000316"
$ws.Range("E11").Value = "  -13.07%  "
Set-TextValue "D12" "42.98"
$ws.Range("E12").Value = "  +4.90%  "
Set-TextValue "D13" "10.42"
$ws.Range("E13").Value = "  +2.84%  "
Set-TextValue "D14" "4.382.29"
$ws.Range("E14").Value = "  -2.44%  "
Set-TextValue "D15" "14.75"
$ws.Range("E15").Value = "  -6.03%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D16" "0.137"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.769.71"
$ws.Range("E17").Value = "  -2.91%  "
Set-TextValue "D18" "19.86"
$ws.Range("E18").Value = "  +0.95%  "
Set-TextValue "D19" "1.13"
$ws.Range("E19").Value = "  +5.36%  "
Set-TextValue "D20" "66.731.01"
$ws.Range("E20").Value = "  -0.71%  "
Set-TextValue "D21" "416.70"
$ws.Range("E21").Value = "  +1.99%  "
Set-TextValue "D22" "14.50"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +7.56%  "
Set-TextValue "D24" "86.02"
$ws.Range("E24").Value = "  +0.81%  "
Set-TextValue "D25" "37.13"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  +5.15%  "
Set-TextValue "D27" "5.64"
$ws.Range("E27").Value = "  -0.48%  "
Set-TextValue "D28" "9.49"
$ws.Range("E28").Value = "  +32.69%  "
Set-TextValue "D29" "9.74"
$ws.Range("E29").Value = "  +1.50%  "
Set-TextValue "D30" "723.52"
$ws.Range("E30").Value = "  +4.78%  "
Set-TextValue "D31" "13.78"
$ws.Range("E31").Value = "  +10.37%  "
$ws.Range("E32").Value = "  +8.88%  "
$ws.Range("E33").Value = "  +2.12%  "
Set-TextValue "D34" "43.29"
$ws.Range("E34").Value = "  +11.50%  "
Set-TextValue "D35" "0.155"
$ws.Range("E35").Value = "  +1.41%  "
Set-TextValue "D36" "0.999"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +24.42%  "
Set-TextValue "D38" "56.49"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D40" "2.91"
$ws.Range("E40").Value = "  -5.22%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D41" "2.65"
$ws.Range("E41").Value = "  +35.83%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D42" "3.36"
$ws.Range("E42").Value = "  +8.34%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D43" "0.141"
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D44" "0.0₃0678"
$ws.Range("E44").Value = "  -16.39%  "
$ws.Range("E45").Value = "  +0.02%  "
Set-TextValue "D46" "0.326"
$ws.Range("E46").Value = "  +11.86%  "
Set-TextValue "D47" "3.30"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("E49").Value = "  +3.92%  "
Set-TextValue "D50" "142.70"
$ws.Range("E50").Value = "  -3.61%  "
Set-TextValue "D51" "2.83"
$ws.Range("E51").Value = "  +1.71%  "
